$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "EMBARGO DATE" column (column C) entirely, shifting all
# subsequent columns one to the left.
$ws.Range("C:C").EntireColumn.Delete()
